$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.1.2"
$meta.Range("B5").Value = "CodeSystem - Transplant Timeline - NMDP"
$meta.Range("B8").Value = "2025-04-15T15:35:56-05:00"

# --- Concepts sheet updates (code/display reorder + renaming) ---
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Range("B2").Value = "12W-6M"
$concepts.Range("C2").Value = "Over 12 weeks - up to 6 months"

$concepts.Range("B3").Value = "4-6W"
$concepts.Range("C3").Value = "Between 4-6 weeks"

$concepts.Range("B4").Value = "4W"
$concepts.Range("C4").Value = "Less than 4 weeks"

$concepts.Range("B5").Value = "6MG"
$concepts.Range("C5").Value = "Greater than 6 months"

$concepts.Range("B6").Value = "7-12w"
$concepts.Range("C6").Value = "Between 7-12 weeks"

$concepts.Range("B7").Value = "NA"
$concepts.Range("C7").Value = "N/A : MUD Transplant not preferred treatment"

$concepts.Range("B8").Value = "PEND"
$concepts.Range("C8").Value = "Pending, Case manager to follow up"
